$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2, 1).Value = "Última actualización: 08:11:27"
$ws.Cells.Item(3, 1).Value = "Total filas: 88"
$ws.Cells.Item(37, 1).Value = "05:42:22"
$ws.Cells.Item(37, 3).Value = "17X38_ROMERO"
$ws.Cells.Item(37, 4).Value = 114
$ws.Cells.Item(38, 1).Value = "06:33:46"
$ws.Cells.Item(38, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(38, 4).Value = 63
$ws.Cells.Item(48, 1).Value = "06:33:46"
$ws.Cells.Item(48, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(48, 4).Value = 87
$ws.Cells.Item(49, 1).Value = "06:16:15"
$ws.Cells.Item(49, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(49, 4).Value = 104
$ws.Cells.Item(57, 1).Value = "08:11:27"
$ws.Cells.Item(57, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(57, 4).Value = 3
$ws.Cells.Item(58, 1).Value = "07:48:35"
$ws.Cells.Item(58, 3).Value = "10_OLMOS"
$ws.Cells.Item(58, 4).Value = 26
$ws.Cells.Item(59, 1).Value = "07:36:59"
$ws.Cells.Item(59, 2).Value = "08:14"
$ws.Cells.Item(59, 4).Value = 38
$ws.Cells.Item(60, 1).Value = "06:16:15"
$ws.Cells.Item(60, 2).Value = "08:15"
$ws.Cells.Item(60, 3).Value = "17_ROMERO"
$ws.Cells.Item(60, 4).Value = 119
$ws.Cells.Item(61, 1).Value = "07:36:59"
$ws.Cells.Item(61, 2).Value = "08:25"
$ws.Cells.Item(61, 4).Value = 49
$ws.Cells.Item(62, 2).Value = "08:26"
$ws.Cells.Item(62, 3).Value = "15X38_ABASTO"
$ws.Cells.Item(62, 4).Value = 113
$ws.Cells.Item(63, 1).Value = "06:33:46"
$ws.Cells.Item(63, 2).Value = "08:27"
$ws.Cells.Item(63, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(63, 4).Value = 114
$ws.Cells.Item(64, 1).Value = "06:45:50"
$ws.Cells.Item(64, 2).Value = "08:29"
$ws.Cells.Item(64, 3).Value = "14_ABASTO"
$ws.Cells.Item(64, 4).Value = 104
$ws.Cells.Item(65, 1).Value = "07:36:59"
$ws.Cells.Item(65, 2).Value = "08:30"
$ws.Cells.Item(65, 4).Value = 54
$ws.Cells.Item(66, 1).Value = "06:33:46"
$ws.Cells.Item(66, 2).Value = "08:31"
$ws.Cells.Item(66, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(66, 4).Value = 118
$ws.Cells.Item(67, 1).Value = "08:11:27"
$ws.Cells.Item(67, 2).Value = "08:33"
$ws.Cells.Item(67, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(67, 4).Value = 22
$ws.Cells.Item(68, 1).Value = "06:45:50"
$ws.Cells.Item(68, 2).Value = "08:38"
$ws.Cells.Item(68, 3).Value = "215C_EL PATO"
$ws.Cells.Item(68, 4).Value = 113
$ws.Cells.Item(69, 2).Value = "08:39"
$ws.Cells.Item(69, 3).Value = "215C_EL PATO"
$ws.Cells.Item(69, 4).Value = 51
$ws.Cells.Item(70, 2).Value = "08:43"
$ws.Cells.Item(70, 3).Value = "10_OLMOS"
$ws.Cells.Item(70, 4).Value = 91
$ws.Cells.Item(71, 2).Value = "08:44"
$ws.Cells.Item(71, 3).Value = "10_OLMOS"
$ws.Cells.Item(71, 4).Value = 56
$ws.Cells.Item(72, 2).Value = "08:49"
$ws.Cells.Item(72, 3).Value = "215A_EL PATO"
$ws.Cells.Item(72, 4).Value = 97
$ws.Cells.Item(73, 1).Value = "07:48:35"
$ws.Cells.Item(73, 2).Value = "08:50"
$ws.Cells.Item(73, 3).Value = "215A_EL PATO"
$ws.Cells.Item(73, 4).Value = 62
$ws.Cells.Item(74, 1).Value = "08:11:27"
$ws.Cells.Item(74, 2).Value = "08:53"
$ws.Cells.Item(74, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(74, 4).Value = 42
$ws.Cells.Item(75, 2).Value = "08:59"
$ws.Cells.Item(75, 3).Value = "215B_EL PATO"
$ws.Cells.Item(75, 4).Value = 107
$ws.Cells.Item(76, 2).Value = "09:01"
$ws.Cells.Item(76, 3).Value = "17X38_ROMERO"
$ws.Cells.Item(76, 4).Value = 85
$ws.Cells.Item(77, 1).Value = "07:36:59"
$ws.Cells.Item(77, 2).Value = "09:02"
$ws.Cells.Item(77, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(77, 4).Value = 86
$ws.Cells.Item(78, 1).Value = "07:12:53"
$ws.Cells.Item(78, 2).Value = "09:02"
$ws.Cells.Item(78, 3).Value = "17X38_ROMERO"
$ws.Cells.Item(78, 4).Value = 110
$ws.Cells.Item(79, 2).Value = "09:04"
$ws.Cells.Item(79, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(79, 4).Value = 88
$ws.Cells.Item(80, 1).Value = "07:48:35"
$ws.Cells.Item(80, 2).Value = "09:08"
$ws.Cells.Item(80, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(80, 4).Value = 80
$ws.Cells.Item(81, 1).Value = "07:36:59"
$ws.Cells.Item(81, 2).Value = "09:14"
$ws.Cells.Item(81, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(81, 4).Value = 98
$ws.Cells.Item(82, 2).Value = "09:14"
$ws.Cells.Item(82, 3).Value = "15_ABASTO"
$ws.Cells.Item(82, 4).Value = 98
$ws.Cells.Item(83, 1).Value = "07:36:59"
$ws.Cells.Item(83, 2).Value = "09:16"
$ws.Cells.Item(83, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(83, 4).Value = 100
$ws.Cells.Item(84, 1).Value = "07:55:46"
$ws.Cells.Item(84, 2).Value = "09:21"
$ws.Cells.Item(84, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(84, 4).Value = 86
$ws.Cells.Item(85, 1).Value = "07:36:59"
$ws.Cells.Item(85, 2).Value = "09:26"
$ws.Cells.Item(85, 3).Value = "215_EL PELIGRO"
$ws.Cells.Item(85, 4).Value = 110
$ws.Cells.Item(86, 2).Value = "09:27"
$ws.Cells.Item(86, 3).Value = "215_EL PELIGRO"
$ws.Cells.Item(86, 4).Value = 99
$ws.Cells.Item(87, 1).Value = "07:36:59"
$ws.Cells.Item(87, 2).Value = "09:30"
$ws.Cells.Item(87, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(87, 4).Value = 114
$ws.Cells.Item(88, 1).Value = "07:48:35"
$ws.Cells.Item(88, 2).Value = "09:31"
$ws.Cells.Item(88, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(88, 4).Value = 103
$ws.Cells.Item(89, 1).Value = "07:48:35"
$ws.Cells.Item(89, 2).Value = "09:39"
$ws.Cells.Item(89, 3).Value = "15_ABASTO"
$ws.Cells.Item(89, 4).Value = 111
$ws.Cells.Item(89, 5).Value = "LP1912"
$ws.Cells.Item(90, 1).Value = "07:48:35"
$ws.Cells.Item(90, 2).Value = "09:44"
$ws.Cells.Item(90, 3).Value = "14_ABASTO"
$ws.Cells.Item(90, 4).Value = 116
$ws.Cells.Item(90, 5).Value = "LP1912"
$ws.Cells.Item(91, 1).Value = "07:55:46"
$ws.Cells.Item(91, 2).Value = "09:51"
$ws.Cells.Item(91, 3).Value = "15_ABASTO"
$ws.Cells.Item(91, 4).Value = 116
$ws.Cells.Item(91, 5).Value = "LP1912"
$ws.Cells.Item(92, 1).Value = "08:11:27"
$ws.Cells.Item(92, 2).Value = "10:03"
$ws.Cells.Item(92, 3).Value = "215C_EL PATO"
$ws.Cells.Item(92, 4).Value = 112
$ws.Cells.Item(92, 5).Value = "LP1912"
$ws.Cells.Item(93, 1).Value = "08:11:27"
$ws.Cells.Item(93, 2).Value = "10:10"
$ws.Cells.Item(93, 3).Value = "10_OLMOS"
$ws.Cells.Item(93, 4).Value = 119
$ws.Cells.Item(93, 5).Value = "LP1912"

$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2, 1).Value = "Última actualización: 08:11:27"
$ws.Cells.Item(3, 1).Value = "Total filas: 16"
$ws.Cells.Item(21, 1).Value = "08:11:27"
$ws.Cells.Item(21, 2).Value = "10:03"
$ws.Cells.Item(21, 3).Value = "215C_EL PATO"
$ws.Cells.Item(21, 4).Value = 112
$ws.Cells.Item(21, 5).Value = "LP1912"

$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2, 1).Value = "Última actualización: 08:11:27"
$ws.Cells.Item(3, 1).Value = "Total filas: 11"
$ws.Cells.Item(9, 1).Value = "08:11:27"
$ws.Cells.Item(9, 2).Value = "08:11"
$ws.Cells.Item(9, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(9, 5).Value = "L6173"
$ws.Cells.Item(10, 1).Value = "06:33:46"
$ws.Cells.Item(10, 2).Value = "08:22"
$ws.Cells.Item(10, 4).Value = 109
$ws.Cells.Item(11, 1).Value = "07:48:35"
$ws.Cells.Item(11, 2).Value = "08:25"
$ws.Cells.Item(11, 4).Value = 37
$ws.Cells.Item(12, 1).Value = "07:55:46"
$ws.Cells.Item(12, 2).Value = "08:26"
$ws.Cells.Item(12, 4).Value = 31
$ws.Cells.Item(13, 2).Value = "08:27"
$ws.Cells.Item(13, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(13, 4).Value = 51
$ws.Cells.Item(13, 5).Value = "L6203"
$ws.Cells.Item(14, 1).Value = "07:36:59"
$ws.Cells.Item(14, 2).Value = "08:51"
$ws.Cells.Item(14, 4).Value = 75
$ws.Cells.Item(15, 1).Value = "07:48:35"
$ws.Cells.Item(15, 2).Value = "08:52"
$ws.Cells.Item(15, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(15, 4).Value = 64
$ws.Cells.Item(15, 5).Value = "L6173"
$ws.Cells.Item(16, 1).Value = "08:11:27"
$ws.Cells.Item(16, 2).Value = "10:09"
$ws.Cells.Item(16, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(16, 4).Value = 118
$ws.Cells.Item(16, 5).Value = "L6203"

Write-Output "done"